$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: standalone run "reminder" -> "remainder"
#   (" (" + "reminder" + ");" are three separate runs that share identical
#    run formatting; a plain text replace would cause the engine to
#    re-coalesce them into a single run, so we briefly diverge the target
#    run's formatting, edit its text, then restore the formatting - a pure
#    property change does not trigger the run-merge pass.)
# ---------------------------------------------------------------------------
$content = $d.Content.Text
$idx = $content.IndexOf("reminder")

$run1 = $d.Range($idx, $idx + 8)
$run1.Font.Bold = 1
$run1b = $d.Range($idx, $idx + 8)
$run1b.Text = "remainder"
$run1c = $d.Range($idx, $idx + 9)
$run1c.Font.Bold = 0

# ---------------------------------------------------------------------------
# Change 2: " - the reminder of an integer division;" (one run) becomes four
# runs: " - the ", "remainder", " ", "of an integer division;"
# ---------------------------------------------------------------------------
$content = $d.Content.Text
$idx2 = $content.IndexOf(" - the reminder of an integer division;")

# Edit the text first (starting one character in, so the edit does not begin
# exactly on the existing run boundary with the preceding ")" run, which
# would otherwise mis-attribute the new text to that preceding run).
$editRng = $d.Range($idx2 + 1, $idx2 + 39)
$editRng.Text = "- the remainder of an integer division;"

# Now split the single run into four via temporary formatting divergence,
# then restore — this creates genuine run boundaries without merging any of
# them back together (or into the neighboring ")" / line-break runs).
$p1 = $d.Range($idx2 + 0, $idx2 + 7)    # " - the "
$p2 = $d.Range($idx2 + 7, $idx2 + 16)   # "remainder"
$p3 = $d.Range($idx2 + 16, $idx2 + 17)  # " "
$p4 = $d.Range($idx2 + 17, $idx2 + 40)  # "of an integer division;"
$pbr = $d.Range($idx2 + 40, $idx2 + 41) # the following line-break run

$p1.Font.Bold = 1
$p2.Font.Bold = 0
$p3.Font.Bold = 1
$p4.Font.Bold = 1
$pbr.Font.Bold = 0

$p1.Font.Bold = 0
$p2.Font.Bold = 0
$p3.Font.Bold = 0
$p4.Font.Bold = 0
